$d = $word.ActiveDocument
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>2024-04-22</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Keep On Learning TDD</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t xml:space="preserve">After submitting my dissertation on the 20th April 2024, I keep on learning TDD. It is evident that learning a new skill is not as easy as I expected, especially by myself. Since TDD is not widely used amongst the software engineers in China, to find out tutorials is more difficult than when you learn a programming language or a Spring framework. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t xml:space="preserve">As a matter of fact, many Chinese developers don't like writing tests, not mention to units test or TDD. Fortunately, I have found a book written by an experienced programmer who had encountered the same obstacles as I have now and he had already overcome them. Additionally, I have found an online course of TDD which is made by the Chief Technology Officer of Thought Works in China. He is an exceptional software engineer and is exactly the person I am looking for. For a intermediate programmer like me, it is impossible to have a mentor teach me advanced skills to be an exceptional developer. Not only the knowledge of TDD have I learned, but also some practical programming approaches such as how to use IDE to refactor effectively, how to write Java code with its new features and so forth. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>In conclusion, though both the book and the course are not easy to understand, I have finished half of the former and two chapters of the latter. Presumably, reading one time of the book is not enough and the course should be watched one or two more times. First of all, finish them for the first time.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
'@
$d.Content.InsertXML($xml)
